$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper values
# ---------------------------------------------------------------------------
$hyperlinkColor = 15570276   # BGR int for FF6495ED (matches existing HyperLink font color)
$dateFmt = "yyyy-mm-dd HH:mm:ss"

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = 2
    $rng.Font.Color = $hyperlinkColor
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
}

function Style-AsDate($rng) {
    $rng.NumberFormat = $dateFmt
}

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1): columns A (File Name, link), B (zh-cn),
# C (de-de), D (Latest Handoff Date)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Push the existing data row (row 2) down to row 3, preserving its styles.
$wsOverview.Rows.Item(2).Insert()

# Remove stale hyperlink definitions (their `ref` still points at the old
# row 2 position after the insert) - we'll recreate them in the right order.
$wsOverview.Hyperlinks.Delete()

# New row 2: newly-handed-off file.
$wsOverview.Range("A2").Value = "2392fc49-9f81-4cb9-8b2a-d566500e7f11.md"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-29-19 18:29:32"
Style-AsHyperlink($wsOverview.Range("A2"))

# Row 3: original file row (values re-asserted defensively after the insert).
$wsOverview.Range("A3").Value = "5f7f7f71-ae82-4c3e-a8fb-26ca246b321e.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-29-19 18:29:12"
Style-AsHyperlink($wsOverview.Range("A3"))

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ac8b5ab35d14987de62ba83aefb1b2afbaa3e6fb/e2e/2392fc49-9f81-4cb9-8b2a-d566500e7f11.md", [System.Type]::Missing, [System.Type]::Missing, "2392fc49-9f81-4cb9-8b2a-d566500e7f11.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ac8b5ab35d14987de62ba83aefb1b2afbaa3e6fb/e2e/5f7f7f71-ae82-4c3e-a8fb-26ca246b321e.md", [System.Type]::Missing, [System.Type]::Missing, "5f7f7f71-ae82-4c3e-a8fb-26ca246b321e.md") | Out-Null

# ---------------------------------------------------------------------------
# Locale detail sheets (zh-cn / de-de): columns A (Source File Name, link),
# B (File Extension, link), C (Status), D (Latest Handoff File, link),
# E (Latest Handoff Datetime), H (Latest Handback DateTime),
# I (Handoff Reason)
# ---------------------------------------------------------------------------
$locales = @(
    @{ SheetName = "zh-cn"; NewXlf = "2392fc49-9f81-4cb9-8b2a-d566500e7f11.05badabfaad0a8f0273f1e4f4aa31b67c0bd129c.zh-cn.xlf"; NewDt = "2016-03-19 18:29:30"; OldXlf = "5f7f7f71-ae82-4c3e-a8fb-26ca246b321e.80bc0d84dfc2f43e600f1410afe7aacc41568e88.zh-cn.xlf"; OldDt = "2016-03-19 18:29:09"; XlfRelFolder = "oltest.zh-cn"; XlfCommit = "1f4f53fd856f6b6c52ce64da6103c6af606dc872" },
    @{ SheetName = "de-de"; NewXlf = "2392fc49-9f81-4cb9-8b2a-d566500e7f11.05badabfaad0a8f0273f1e4f4aa31b67c0bd129c.de-de.xlf"; NewDt = "2016-03-19 18:29:32"; OldXlf = "5f7f7f71-ae82-4c3e-a8fb-26ca246b321e.80bc0d84dfc2f43e600f1410afe7aacc41568e88.de-de.xlf"; OldDt = "2016-03-19 18:29:12"; XlfRelFolder = "oltest.de-de"; XlfCommit = "2b1bc1a86e0e657dd82c216505a725875e09bdc4" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.SheetName)

    # Push the existing data row (row 2) down to row 3.
    $ws.Rows.Item(2).Insert()

    # Drop stale hyperlink definitions; rebuilt below in correct order.
    $ws.Hyperlinks.Delete()

    # New row 2: newly handed-off file.
    $ws.Range("A2").Value = "2392fc49-9f81-4cb9-8b2a-d566500e7f11.md"
    $ws.Range("B2").Value = ".md"
    $ws.Range("C2").Value = "Ready for handoff"
    $ws.Range("D2").Value = $loc.NewXlf
    $ws.Range("E2").Value = $loc.NewDt
    $ws.Range("H2").Value = "0001-01-01 00:00:00"
    $ws.Range("I2").Value = "Include"
    Style-AsHyperlink($ws.Range("A2"))
    Style-AsHyperlink($ws.Range("B2"))
    Style-AsHyperlink($ws.Range("D2"))
    Style-AsDate($ws.Range("E2"))

    # Row 3: original file row (values re-asserted defensively).
    $ws.Range("A3").Value = "5f7f7f71-ae82-4c3e-a8fb-26ca246b321e.md"
    $ws.Range("B3").Value = ".md"
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Range("D3").Value = $loc.OldXlf
    $ws.Range("E3").Value = $loc.OldDt
    $ws.Range("H3").Value = "0001-01-01 00:00:00"
    $ws.Range("I3").Value = "Include"
    Style-AsHyperlink($ws.Range("A3"))
    Style-AsHyperlink($ws.Range("B3"))
    Style-AsHyperlink($ws.Range("D3"))
    Style-AsDate($ws.Range("E3"))

    $mdUrlNew = "https://github.com/OpenLocalizationTest/oltest/blob/ac8b5ab35d14987de62ba83aefb1b2afbaa3e6fb/e2e/2392fc49-9f81-4cb9-8b2a-d566500e7f11.md"
    $mdUrlOld = "https://github.com/OpenLocalizationTest/oltest/blob/ac8b5ab35d14987de62ba83aefb1b2afbaa3e6fb/e2e/5f7f7f71-ae82-4c3e-a8fb-26ca246b321e.md"
    $xlfUrlNew = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $loc.XlfCommit + "/ol-handoff/OpenLocalizationTestOrg/" + $loc.XlfRelFolder + "/ci/ht/" + $loc.NewXlf
    $xlfUrlOld = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $loc.XlfCommit + "/ol-handoff/OpenLocalizationTestOrg/" + $loc.XlfRelFolder + "/ci/ht/" + $loc.OldXlf

    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrlNew, [System.Type]::Missing, [System.Type]::Missing, "2392fc49-9f81-4cb9-8b2a-d566500e7f11.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B2"), $mdUrlNew, [System.Type]::Missing, [System.Type]::Missing, ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D2"), $xlfUrlNew, [System.Type]::Missing, [System.Type]::Missing, $loc.NewXlf) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A3"), $mdUrlOld, [System.Type]::Missing, [System.Type]::Missing, "5f7f7f71-ae82-4c3e-a8fb-26ca246b321e.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B3"), $mdUrlOld, [System.Type]::Missing, [System.Type]::Missing, ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D3"), $xlfUrlOld, [System.Type]::Missing, [System.Type]::Missing, $loc.OldXlf) | Out-Null
}

"Report regenerated for handoff"
